# Apply cell-value changes to Sheet1 per the target diff (state-machine table update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "AvvioChiusuraPostSelfCare"
$ws.Range("C10").Value = "AvvioChiusura"
$ws.Range("B11").Clear()
$ws.Range("C11").Value = "AvvioPerizia"
$ws.Range("B14").Value = "VisioCompletataNuovaVisio"
$ws.Range("C14").Value = "AvvioManualeVisio"
$ws.Range("B15").Value = "AppuntamentoModificato"
$ws.Range("C15").Value = "EsecuzioneVisio"
$ws.Range("B16").Value = "AppuntamentoAnnullato"
$ws.Range("B17").Value = "DatiObbligatoriMancanti"
$ws.Range("B18").Value = "VisioRifiutata"
$ws.Range("C18").Value = "RiassegnazioneVisio"
$ws.Range("B19").Value = "VisioCompletata"
$ws.Range("C19").Value = "AvvioPerizia"
$ws.Range("A21").Value = "Pa_Visio"
$ws.Range("A22").Value = "Pa_Visio"
$ws.Range("B25").Value = "DeskCompletataNuovaDesk"
$ws.Range("C25").Value = "AvvioManualeDesk"
$ws.Range("B26").Value = "AppuntamentoModificato"
$ws.Range("C26").Value = "EsecuzioneDesk"
$ws.Range("B27").Value = "AppuntamentoAnnullato"
$ws.Range("B28").Value = "DatiObbligatoriMancanti"
$ws.Range("B29").Value = "DeskRifiutata"
$ws.Range("C29").Value = "RiassegnazioneDesk"
$ws.Range("B30").Value = "DeskCompletata"
$ws.Range("C30").Value = "AvvioPeriziaPostDesk"
$ws.Range("A32").Value = "Pa_Desk"
$ws.Range("C32").Value = "EsecuzioneDesk"
$ws.Range("A33").Value = "Pa_Desk"
$ws.Range("C33").Value = "EsecuzioneDesk"
$ws.Range("B38").Value = "AnnullamentoSopralluogo"
$ws.Range("C38").Value = "AssegnazioneSopralluogo"
$ws.Range("B39").Value = "AppuntamentoAnnullato"
$ws.Range("C39").Value = "PianificazioneSopralluogo"
$ws.Range("B40").Value = "ChiusuraAttivita"
$ws.Range("B41").Value = "ChiusuraSoprallVariato"
$ws.Range("B42").Value = "AppuntamentoModificato"
$ws.Range("C42").Value = "EsecuzioneSopralluogo"
$ws.Range("B43").Value = "DatiObbligatoriMancanti"
$ws.Range("B44").Value = "RichiestaNuovoSopralluogo"
$ws.Range("C44").Value = "ApprovazioneSopralluogo"
$ws.Range("B45").Value = "AnnullamentoSopralluogo"
$ws.Range("C45").Value = "AssegnazioneSopralluogoPCE"
$ws.Range("C47").Value = "PianificazioneSopralluogo"
$ws.Range("A49").Value = "Pa_Sopralluogo"
$ws.Range("C49").Value = "EsecuzioneSopralluogo"
$ws.Range("A50").Value = "Pa_Sopralluogo"
$ws.Range("C50").Value = "EsecuzioneSopralluogo"
$ws.Range("A51").Value = "AvvioSopralluogoPCE"
$ws.Range("B51").Value = "SopralluogoAvvioto"
$ws.Range("C51").Value = "AssegnazioneSopralluogoPCE"
$ws.Range("A52").Value = "AssegnazioneSopralluogoPCE"
$ws.Range("B52").Value = "SopralluogoAssegnato"
$ws.Range("C52").Value = "PianificazioneSopralluogoPCE"
$ws.Range("A53").Value = "PianificazioneSopralluogoPCE"
$ws.Range("B53").Value = "SopralluogoPianificato"
$ws.Range("C53").Value = "EsecuzioneSopralluogoPCE"
$ws.Range("B54").Value = "SopralluogoRifiutato"
$ws.Range("C54").Value = "VerificaSopralluogoPCEDopoRifiuto"
$ws.Range("A55").Value = "EsecuzioneSopralluogoPCE"
$ws.Range("B55").Value = "AppuntamentoAnnullato"
$ws.Range("C55").Value = "PianificazioneSopralluogoPCE"
$ws.Range("B56").Value = "ChiusuraAttivita"
$ws.Range("C56").Value = "PianificazioneSopralluogoPCE"
$ws.Range("B57").Value = "ChiusuraSoprallVariato"
$ws.Range("C57").Value = "PianificazioneSopralluogoPCE"
$ws.Range("B58").Value = "AppuntamentoModificato"
$ws.Range("C58").Value = "EsecuzioneSopralluogoPCE"
$ws.Range("B59").Value = "DatiObbligatoriMancanti"
$ws.Range("C59").Value = "EsecuzioneSopralluogoPCE"
$ws.Range("B60").Value = "RichiestaNuovoSopralluogo"
$ws.Range("C60").Value = "ApprovazioneSopralluogoPCE"
$ws.Range("A61").Value = "VerificaSopralluogoPCEDopoRifiuto"
$ws.Range("B61").Value = "SopralluogoAssegnato"
$ws.Range("C61").Value = "PianificazioneSopralluogoPCE"
$ws.Range("A62").Value = "ApprovazioneSopralluogoPCE"
$ws.Range("B62").Value = "Approvato"
$ws.Range("C62").Value = "PianificazioneSopralluogoPCE"
$ws.Range("A63").Value = "ApprovazioneSopralluogoPCE"
$ws.Range("B63").Value = "NonApprovato"
$ws.Range("C63").Value = "NuovoSopralluogoNonApprovatoPCE"
$ws.Range("A64").Value = "Pa_SopralluogoPCE"
$ws.Range("B64").Value = "RichiestaPaRicevuta"
$ws.Range("C64").Value = "EsecuzioneSopralluogoPCE"
$ws.Range("A65").Value = "Pa_SopralluogoPCE"
$ws.Range("B65").Value = "ConfermaPaRicevuta"
$ws.Range("C65").Value = "EsecuzioneSopralluogoPCE"

# Remove rows 66-68 entirely (no longer present in the target state machine)
$ws.Rows("66:68").Delete()

Write-Host "assessment sheet updated"
